$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cells, written in the same order the original authoring tool
# appended them to the shared-string table (H1, G2, H2, G1)
$ws.Range("H1").Value = "Comment"
$ws.Range("G2").Value = "CM - Cost Structure"
$ws.Range("H2").Value = "cost structure"
$ws.Range("G1").Value = "Actions"

# Copy style from existing header/data cells so new cells match formatting
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("F2").Copy()
$ws.Range("G2:H2").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Update the selection to match the target state
$ws.Range("E14").Select()
